{"js": "// Fixed disease contrast table label (2014-01 to 2022-12)\n\n// 1. Update every \"Trend difference\" table-label cell from the\n//    2012-01/2021-12 date range to 2014-01/2022-12 across the document.\nconst results = context.document.body.search(\n  \"Trend difference (2012-01 to 2021-12)\",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\n    \"Trend difference (2014-01 to 2022-12)\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 2. Tighten the space-before on the existing \"Abstract\" paragraph style\n//    (300 -> 100 twips, i.e. 15pt -> 5pt).\nconst existingStyles = context.document.getStyles();\nconst abstractStyle = existingStyles.getByNameOrNullObject(\"Abstract\");\nawait context.sync();\nabstractStyle.paragraphFormat.spaceBefore = 5;\nawait context.sync();\n\n// 3. Add the new \"Abstract Title\" paragraph style.\ncontext.document.addStyle(\"Abstract Title\", \"Paragraph\");\nawait context.sync();\n\nconst abstractTitle = context.document\n  .getStyles()\n  .getByNameOrNullObject(\"Abstract Title\");\nawait context.sync();\n\nabstractTitle.baseStyle = \"Normal\";\nabstractTitle.nextParagraphStyle = \"Abstract\";\nabstractTitle.quickStyle = true;\nabstractTitle.paragraphFormat.keepWithNext = true;\nabstractTitle.paragraphFormat.keepTogether = true;\nabstractTitle.paragraphFormat.alignment = \"Centered\";\nabstractTitle.paragraphFormat.spaceBefore = 15;\nabstractTitle.paragraphFormat.spaceAfter = 0;\nabstractTitle.font.size = 10;\nabstractTitle.font.bold = true;\nabstractTitle.font.color = \"#345A8A\";\nawait context.sync();\n\n// 4. Add the new \"Footnote Block Text\" paragraph style.\ncontext.document.addStyle(\"Footnote Block Text\", \"Paragraph\");\nawait context.sync();\n\nconst footnoteBlockText = context.document\n  .getStyles()\n  .getByNameOrNullObject(\"Footnote Block Text\");\nawait context.sync();\n\nfootnoteBlockText.baseStyle = \"Footnote Text\";\nfootnoteBlockText.nextParagraphStyle = \"Footnote Text\";\nfootnoteBlockText.priority = 9;\nfootnoteBlockText.unhideWhenUsed = true;\nfootnoteBlockText.quickStyle = true;\nfootnoteBlockText.paragraphFormat.spaceBefore = 5;\nfootnoteBlockText.paragraphFormat.spaceAfter = 5;\nfootnoteBlockText.paragraphFormat.firstLineIndent = 0;\nfootnoteBlockText.paragraphFormat.leftIndent = 24;\nfootnoteBlockText.paragraphFormat.rightIndent = 24;\nawait context.sync();\n", "ps1": "# Fixed disease contrast table label (2014-01 to 2022-12)\n$d = $word.ActiveDocument\n\n# 1. Update every \"Trend difference\" table-label cell from the\n#    2012-01/2021-12 date range to 2014-01/2022-12 across the document.\n$find = $d.Content.Find\n$find.Execute(\n    \"Trend difference (2012-01 to 2021-12)\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"Trend difference (2014-01 to 2022-12)\",\n    2\n) | Out-Null\n\n# 2. Tighten the space-before on the existing \"Abstract\" paragraph style\n#    (300 -> 100 twips, i.e. 15pt -> 5pt).\n$abstractStyle = $d.Styles(\"Abstract\")\n$abstractStyle.ParagraphFormat.SpaceBefore = 5\n\n# 3. Add the new \"Abstract Title\" paragraph style.\n$abstractTitle = $d.Styles.Add(\"AbstractTitle\", 1)\n$abstractTitle.NameLocal = \"Abstract Title\"\n$abstractTitle.BaseStyle = $d.Styles(\"Normal\")\n$abstractTitle.NextParagraphStyle = $d.Styles(\"Abstract\")\n$abstractTitle.QuickStyle = $true\n$abstractTitle.ParagraphFormat.KeepWithNext = $true\n$abstractTitle.ParagraphFormat.KeepTogether = $true\n$abstractTitle.ParagraphFormat.Alignment = 1\n$abstractTitle.ParagraphFormat.SpaceAfter = 0\n$abstractTitle.ParagraphFormat.SpaceBefore = 15\n$abstractTitle.Font.Size = 10\n$abstractTitle.Font.SizeBi = 10\n$abstractTitle.Font.Bold = $true\n$abstractTitle.Font.Color = 9067060\n\n# 4. Add the new \"Footnote Block Text\" paragraph style.\n$footnoteBlockText = $d.Styles.Add(\"FootnoteBlockText\", 1)\n$footnoteBlockText.NameLocal = \"Footnote Block Text\"\n$footnoteBlockText.BaseStyle = $d.Styles(\"Footnote Text\")\n$footnoteBlockText.NextParagraphStyle = $d.Styles(\"Footnote Text\")\n$footnoteBlockText.Priority = 9\n$footnoteBlockText.UnhideWhenUsed = $true\n$footnoteBlockText.QuickStyle = $true\n$footnoteBlockText.ParagraphFormat.SpaceAfter = 5\n$footnoteBlockText.ParagraphFormat.SpaceBefore = 5\n$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0\n$footnoteBlockText.ParagraphFormat.LeftIndent = 24\n$footnoteBlockText.ParagraphFormat.RightIndent = 24\n\nWrite-Output \"done\"\n"}
